$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update recalculated odds values in columns AS (PPG_Home) / AT (PPG_Away) ---
# Global find/replace: every occurrence of the old value in that column becomes the new value
$ws.Cells.Item(2, 45).Value = 2.67
$ws.Cells.Item(12, 46).Value = 0.78
$ws.Cells.Item(15, 45).Value = 2
$ws.Cells.Item(16, 46).Value = 1.11
$ws.Cells.Item(17, 46).Value = 1.56
$ws.Cells.Item(22, 45).Value = 1.78
$ws.Cells.Item(26, 45).Value = 2.67
$ws.Cells.Item(29, 45).Value = 2
$ws.Cells.Item(33, 46).Value = 1.11
$ws.Cells.Item(36, 45).Value = 1.78
$ws.Cells.Item(42, 45).Value = 2.67
$ws.Cells.Item(48, 46).Value = 1.56
$ws.Cells.Item(49, 46).Value = 1.11
$ws.Cells.Item(51, 45).Value = 2
$ws.Cells.Item(57, 45).Value = 1.78
$ws.Cells.Item(58, 45).Value = 2.67
$ws.Cells.Item(68, 46).Value = 1.56
$ws.Cells.Item(69, 46).Value = 0.78
$ws.Cells.Item(70, 46).Value = 1.11
$ws.Cells.Item(77, 45).Value = 2.67
$ws.Cells.Item(81, 45).Value = 2
$ws.Cells.Item(82, 45).Value = 1.78
$ws.Cells.Item(88, 46).Value = 0.78
$ws.Cells.Item(92, 46).Value = 1.11
$ws.Cells.Item(97, 45).Value = 1.78
$ws.Cells.Item(100, 45).Value = 2.67
$ws.Cells.Item(101, 45).Value = 2
$ws.Cells.Item(103, 46).Value = 1.56
$ws.Cells.Item(104, 46).Value = 0.78
$ws.Cells.Item(106, 46).Value = 1.11
$ws.Cells.Item(114, 45).Value = 2.67
$ws.Cells.Item(115, 45).Value = 2
$ws.Cells.Item(116, 45).Value = 1.78
$ws.Cells.Item(118, 46).Value = 1.56
$ws.Cells.Item(123, 46).Value = 0.78
$ws.Cells.Item(124, 45).Value = 1.78
$ws.Cells.Item(127, 46).Value = 1.56
$ws.Cells.Item(129, 45).Value = 2
$ws.Cells.Item(129, 46).Value = 1.11
$ws.Cells.Item(130, 46).Value = 0.78
$ws.Cells.Item(135, 46).Value = 1.56
$ws.Cells.Item(137, 46).Value = 0.78
$ws.Cells.Item(147, 45).Value = 2
$ws.Cells.Item(148, 45).Value = 2.67
$ws.Cells.Item(150, 46).Value = 1.11
$ws.Cells.Item(151, 45).Value = 1.78
$ws.Cells.Item(151, 46).Value = 1.56
$ws.Cells.Item(154, 46).Value = 0.78

# --- 2) Append 3 new match rows (170-172) ---
# Row 170
$ws.Cells.Item(170, 1).Value = 169
$ws.Cells.Item(170, 2).Value = 5053111
$ws.Cells.Item(170, 3).Value = 'Peru Primera División'
$ws.Cells.Item(170, 4).Value = '2023'
$ws.Cells.Item(170, 5).Value = 45088.625
$ws.Cells.Item(170, 6).Value = 19
$ws.Cells.Item(170, 7).Value = 'Alianza Atlético'
$ws.Cells.Item(170, 8).Value = 'Sporting Cristal'
$ws.Cells.Item(170, 9).Value = 0
$ws.Cells.Item(170, 10).Value = 0
$ws.Cells.Item(170, 11).Value = 0
$ws.Cells.Item(170, 12).Value = 0
$ws.Cells.Item(170, 13).Value = 0
$ws.Cells.Item(170, 14).Value = 0
$ws.Cells.Item(170, 15).Value = '[]'
$ws.Cells.Item(170, 16).Value = '[]'
$ws.Cells.Item(170, 17).Value = 1
$ws.Cells.Item(170, 18).Value = 8
$ws.Cells.Item(170, 19).Value = 9
$ws.Cells.Item(170, 20).Value = 4
$ws.Cells.Item(170, 21).Value = 2.25
$ws.Cells.Item(170, 22).Value = 2.6
$ws.Cells.Item(170, 23).Value = 1.35
$ws.Cells.Item(170, 24).Value = 3.2
$ws.Cells.Item(170, 25).Value = 2.6
$ws.Cells.Item(170, 26).Value = 1.43
$ws.Cells.Item(170, 27).Value = 6.5
$ws.Cells.Item(170, 28).Value = 1.1
$ws.Cells.Item(170, 29).Value = 3.55
$ws.Cells.Item(170, 30).Value = 3.6
$ws.Cells.Item(170, 31).Value = 1.97
$ws.Cells.Item(170, 32).Value = 1.03
$ws.Cells.Item(170, 33).Value = 13
$ws.Cells.Item(170, 34).Value = 1.23
$ws.Cells.Item(170, 35).Value = 3.75
$ws.Cells.Item(170, 36).Value = 1.67
$ws.Cells.Item(170, 37).Value = 2
$ws.Cells.Item(170, 38).Value = 1.7
$ws.Cells.Item(170, 39).Value = 2.15
$ws.Cells.Item(170, 40).Value = 1.87
$ws.Cells.Item(170, 41).Value = 1.22
$ws.Cells.Item(170, 42).Value = 1.28
$ws.Cells.Item(170, 43).Value = 2.13
$ws.Cells.Item(170, 44).Value = 1.63
$ws.Cells.Item(170, 45).Value = 2
$ws.Cells.Item(170, 46).Value = 1.56
$ws.Cells.Item(170, 47).Value = 1.57
$ws.Cells.Item(170, 48).Value = 1.36
$ws.Cells.Item(170, 49).Value = 2.93
$ws.Cells.Item(170, 50).Value = 2.63
$ws.Cells.Item(170, 51).Value = 8.5
$ws.Cells.Item(170, 52).Value = 1.62
$ws.Cells.Item(170, 53).Value = 1.2
$ws.Cells.Item(170, 54).Value = 1.32
$ws.Cells.Item(170, 55).Value = 1.56
$ws.Cells.Item(170, 56).Value = 1.9
$ws.Cells.Item(170, 57).Value = 2.62
$ws.Cells.Item(170, 58).Value = 2
$ws.Cells.Item(170, 59).Value = 8
$ws.Cells.Item(170, 60).Value = 5
$ws.Cells.Item(170, 61).Value = 9
$ws.Cells.Item(170, 62).Value = 7
$ws.Cells.Item(170, 63).Value = 17
# Row 171
$ws.Cells.Item(171, 1).Value = 170
$ws.Cells.Item(171, 2).Value = 5053109
$ws.Cells.Item(171, 3).Value = 'Peru Primera División'
$ws.Cells.Item(171, 4).Value = '2023'
$ws.Cells.Item(171, 5).Value = 45088.72916666666
$ws.Cells.Item(171, 6).Value = 19
$ws.Cells.Item(171, 7).Value = 'Real Garcilaso'
$ws.Cells.Item(171, 8).Value = 'César Vallejo'
$ws.Cells.Item(171, 9).Value = 0
$ws.Cells.Item(171, 10).Value = 0
$ws.Cells.Item(171, 11).Value = 0
$ws.Cells.Item(171, 12).Value = 2
$ws.Cells.Item(171, 13).Value = 0
$ws.Cells.Item(171, 14).Value = 2
$ws.Cells.Item(171, 15).Value = '[''54'', ''57'']'
$ws.Cells.Item(171, 16).Value = '[]'
$ws.Cells.Item(171, 17).Value = 7
$ws.Cells.Item(171, 18).Value = 3
$ws.Cells.Item(171, 19).Value = 10
$ws.Cells.Item(171, 20).Value = 2.1
$ws.Cells.Item(171, 21).Value = 2.25
$ws.Cells.Item(171, 22).Value = 6
$ws.Cells.Item(171, 23).Value = 1.4
$ws.Cells.Item(171, 24).Value = 3.1
$ws.Cells.Item(171, 25).Value = 2.63
$ws.Cells.Item(171, 26).Value = 1.4
$ws.Cells.Item(171, 27).Value = 6.5
$ws.Cells.Item(171, 28).Value = 1.08
$ws.Cells.Item(171, 29).Value = 1.53
$ws.Cells.Item(171, 30).Value = 3.7
$ws.Cells.Item(171, 31).Value = 5.25
$ws.Cells.Item(171, 32).Value = 1.03
$ws.Cells.Item(171, 33).Value = 10
$ws.Cells.Item(171, 34).Value = 1.29
$ws.Cells.Item(171, 35).Value = 3.24
$ws.Cells.Item(171, 36).Value = 1.85
$ws.Cells.Item(171, 37).Value = 1.95
$ws.Cells.Item(171, 38).Value = 1.95
$ws.Cells.Item(171, 39).Value = 1.8
$ws.Cells.Item(171, 40).Value = 1.11
$ws.Cells.Item(171, 41).Value = 1.18
$ws.Cells.Item(171, 42).Value = 2.45
$ws.Cells.Item(171, 43).Value = 2.63
$ws.Cells.Item(171, 44).Value = 0.88
$ws.Cells.Item(171, 45).Value = 2.67
$ws.Cells.Item(171, 46).Value = 0.78
$ws.Cells.Item(171, 47).Value = 1.62
$ws.Cells.Item(171, 48).Value = 1.19
$ws.Cells.Item(171, 49).Value = 2.81
$ws.Cells.Item(171, 50).Value = 1.33
$ws.Cells.Item(171, 51).Value = 9.5
$ws.Cells.Item(171, 52).Value = 3.74
$ws.Cells.Item(171, 53).Value = 1.17
$ws.Cells.Item(171, 54).Value = 1.28
$ws.Cells.Item(171, 55).Value = 1.51
$ws.Cells.Item(171, 56).Value = 1.95
$ws.Cells.Item(171, 57).Value = 2.45
$ws.Cells.Item(171, 58).Value = 9
$ws.Cells.Item(171, 59).Value = 0
$ws.Cells.Item(171, 60).Value = 5
$ws.Cells.Item(171, 61).Value = 4
$ws.Cells.Item(171, 62).Value = 14
$ws.Cells.Item(171, 63).Value = 4
# Row 172
$ws.Cells.Item(172, 1).Value = 171
$ws.Cells.Item(172, 2).Value = 5053112
$ws.Cells.Item(172, 3).Value = 'Peru Primera División'
$ws.Cells.Item(172, 4).Value = '2023'
$ws.Cells.Item(172, 5).Value = 45088.83333333334
$ws.Cells.Item(172, 6).Value = 19
$ws.Cells.Item(172, 7).Value = 'Sport Huancayo'
$ws.Cells.Item(172, 8).Value = 'Universitario'
$ws.Cells.Item(172, 9).Value = 0
$ws.Cells.Item(172, 10).Value = 0
$ws.Cells.Item(172, 11).Value = 0
$ws.Cells.Item(172, 12).Value = 1
$ws.Cells.Item(172, 13).Value = 0
$ws.Cells.Item(172, 14).Value = 1
$ws.Cells.Item(172, 15).Value = '[''77'']'
$ws.Cells.Item(172, 16).Value = '[]'
$ws.Cells.Item(172, 17).Value = 2
$ws.Cells.Item(172, 18).Value = 1
$ws.Cells.Item(172, 19).Value = 3
$ws.Cells.Item(172, 20).Value = 3
$ws.Cells.Item(172, 21).Value = 2.2
$ws.Cells.Item(172, 22).Value = 3.4
$ws.Cells.Item(172, 23).Value = 1.4
$ws.Cells.Item(172, 24).Value = 3
$ws.Cells.Item(172, 25).Value = 2.8
$ws.Cells.Item(172, 26).Value = 1.4
$ws.Cells.Item(172, 27).Value = 7
$ws.Cells.Item(172, 28).Value = 1.08
$ws.Cells.Item(172, 29).Value = 2.44
$ws.Cells.Item(172, 30).Value = 3.2
$ws.Cells.Item(172, 31).Value = 2.58
$ws.Cells.Item(172, 32).Value = 1.04
$ws.Cells.Item(172, 33).Value = 12
$ws.Cells.Item(172, 34).Value = 1.3
$ws.Cells.Item(172, 35).Value = 3.6
$ws.Cells.Item(172, 36).Value = 1.81
$ws.Cells.Item(172, 37).Value = 1.9
$ws.Cells.Item(172, 38).Value = 1.72
$ws.Cells.Item(172, 39).Value = 2.05
$ws.Cells.Item(172, 40).Value = 1.38
$ws.Cells.Item(172, 41).Value = 1.25
$ws.Cells.Item(172, 42).Value = 1.65
$ws.Cells.Item(172, 43).Value = 1.63
$ws.Cells.Item(172, 44).Value = 1.25
$ws.Cells.Item(172, 45).Value = 1.78
$ws.Cells.Item(172, 46).Value = 1.11
$ws.Cells.Item(172, 47).Value = 2.09
$ws.Cells.Item(172, 48).Value = 1.33
$ws.Cells.Item(172, 49).Value = 3.42
$ws.Cells.Item(172, 50).Value = 2.1
$ws.Cells.Item(172, 51).Value = 8
$ws.Cells.Item(172, 52).Value = 1.95
$ws.Cells.Item(172, 53).Value = 1.18
$ws.Cells.Item(172, 54).Value = 1.35
$ws.Cells.Item(172, 55).Value = 1.59
$ws.Cells.Item(172, 56).Value = 2
$ws.Cells.Item(172, 57).Value = 2.5
$ws.Cells.Item(172, 58).Value = 6
$ws.Cells.Item(172, 59).Value = 4
$ws.Cells.Item(172, 60).Value = 7
$ws.Cells.Item(172, 61).Value = 7
$ws.Cells.Item(172, 62).Value = 13
$ws.Cells.Item(172, 63).Value = 11

# --- 3) Apply styles to new rows: column A (bold/border/center-top) and column E (date format) ---
$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Cells.Item(170, 1).PasteSpecial(-4122)
$ws.Cells.Item(171, 1).PasteSpecial(-4122)
$ws.Cells.Item(172, 1).PasteSpecial(-4122)
$ws.Cells.Item(2, 5).Copy() | Out-Null
$ws.Cells.Item(170, 5).PasteSpecial(-4122)
$ws.Cells.Item(171, 5).PasteSpecial(-4122)
$ws.Cells.Item(172, 5).PasteSpecial(-4122)

Write-Output "edit complete"